# Chiffres COVID-19 Valais.xlsx - daily data refresh
# Updates "Nb nouveaux cas positifs" (col C) for a few recent days and
# fills in the previously-empty row for 2020-09-01 (row 189) with that
# day's figures, then nudges the view (scroll position / active cell)
# the way the original author's Excel session left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 181 (2020-08-24): new positive cases 15 -> 16 ---
$ws.Range("C181").Value = 16

# --- Row 182 (2020-08-25): new positive cases 14 -> 13 ---
$ws.Range("C182").Value = 13

# --- Row 187 (2020-08-30): new positive cases 1 -> 0 ---
$ws.Range("C187").Value = 0

# --- Row 188 (2020-08-31): new positive cases 0 -> 3 ---
$ws.Range("C188").Value = 3

# --- Row 189 (2020-09-01): was a blank placeholder row, now filled in ---
$ws.Range("D189").Value = 0
$ws.Range("E189").Value = 0
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 3
$ws.Range("I189").Value = 0
$ws.Range("L189").Value = "0"
$ws.Range("M189").Value = "0"

# --- View state: scroll the frozen pane down a bit and move the active
#     cell from N190 to O190, matching where the author left the sheet ---
$excel.ActiveWindow.ScrollRow = 174
$excel.ActiveWindow.ScrollColumn = 2
$null = $ws.Range("O190").Select()
